# Update the "Pais" (countries) worksheet with refreshed case-count data.
# The source data was re-sorted by "Casos totales" (col B) descending; several
# countries changed rank (Israel overtook Irlanda, Timor Oriental / Fiyi / Dominica /
# Montserrat moved up), and the "last updated" timestamp advanced from 08:52 to 09:22.
# Rather than re-implement a generic sort, we apply the exact resulting cell values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is @(row, column, newValue) using the same column numbering as the sheet:
# 1=Pais 2=Casos totales 3=Nuevos casos 4=Casos activos 5=Recuperados 6=Casos criticos 7=Muertes hoy 8=Muertes
$changes = @(
    @(1, 1, "Datos actualizados a 16 de Abril de 2020 a las 09:22"),
    @(21, 1, "Israel"),
    @(21, 2, 12591),
    @(21, 3, 90),
    @(21, 4, 2624),
    @(21, 5, 9827),
    @(21, 6, 174),
    @(21, 7, 10),
    @(21, 8, 140),
    @(22, 1, "Irlanda"),
    @(22, 2, 12547),
    @(22, 4, 77),
    @(22, 5, 12026),
    @(22, 6, 158),
    @(22, 7, 0),
    @(22, 8, 444),
    @(30, 4, 774),
    @(30, 5, 6522),
    @(36, 2, 6303),
    @(36, 3, 2),
    @(36, 4, 831),
    @(36, 5, 5306),
    @(36, 6, 75),
    @(177, 1, "Timor Oriental"),
    @(177, 3, 10),
    @(177, 6, 0),
    @(178, 1, "Nueva Caledonia"),
    @(178, 4, 1),
    @(178, 5, 17),
    @(178, 8, 0),
    @(179, 1, "Belice"),
    @(179, 2, 18),
    @(179, 5, 16),
    @(179, 6, 1),
    @(179, 8, 2),
    @(181, 1, "Islas Virgenes de los Estados Unidos"),
    @(181, 2, 17),
    @(181, 4, 0),
    @(181, 5, 17),
    @(182, 1, "Nepal"),
    @(182, 4, 1),
    @(182, 5, 15),
    @(182, 6, 0),
    @(182, 8, 0),
    @(183, 1, "Malaui"),
    @(183, 4, 0),
    @(183, 5, 14),
    @(183, 6, 1),
    @(183, 8, 2),
    @(184, 1, "Namibia"),
    @(184, 3, 0),
    @(184, 4, 3),
    @(184, 5, 13),
    @(186, 1, "Suazilandia"),
    @(186, 2, 16),
    @(186, 3, 1),
    @(186, 4, 8),
    @(186, 5, 8),
    @(187, 1, "Santa Lucia"),
    @(187, 2, 15),
    @(187, 4, 11),
    @(187, 5, 4),
    @(187, 6, 0),
    @(188, 1, "Granada"),
    @(188, 6, 2),
    @(189, 1, "San Cristobal y Nieves"),
    @(189, 4, 0),
    @(189, 5, 14),
    @(189, 8, 0),
    @(190, 1, "Curazao"),
    @(190, 2, 14),
    @(190, 4, 10),
    @(190, 5, 3),
    @(190, 8, 1),
    @(191, 1, "Sierra Leona"),
    @(191, 5, 13),
    @(191, 8, 0),
    @(192, 1, "Botsuana"),
    @(192, 2, 13),
    @(192, 4, 0),
    @(192, 5, 12),
    @(192, 8, 1),
    @(193, 1, "San Vicente y las Granadinas"),
    @(193, 4, 1),
    @(193, 5, 11),
    @(194, 1, "Republica de Africa Central"),
    @(194, 2, 12),
    @(194, 4, 4),
    @(194, 5, 8),
    @(195, 1, "Seychelles"),
    @(195, 4, 0),
    @(195, 5, 11),
    @(197, 1, "Islas Malvinas"),
    @(197, 4, 1),
    @(197, 5, 10),
    @(198, 1, "Groenlandia"),
    @(198, 2, 11),
    @(198, 4, 11),
    @(198, 5, 0),
    @(198, 8, 0),
    @(199, 1, "Islas Turcas y Caicos"),
    @(199, 4, 0),
    @(199, 5, 9),
    @(200, 1, "Surinam"),
    @(200, 2, 10),
    @(200, 4, 6),
    @(200, 5, 3),
    @(201, 1, "Gambia"),
    @(201, 4, 2),
    @(201, 5, 6),
    @(202, 1, "Nicaragua"),
    @(202, 2, 9),
    @(202, 4, 4),
    @(202, 5, 4),
    @(202, 8, 1)
)

foreach ($chg in $changes) {
    $ws.Cells.Item($chg[0], $chg[1]).Value = $chg[2]
}
